$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.982.27"
$ws.Range("E2").Value = "  -5.33%  "
$ws.Range("D3").Value = "2.218.71"
$ws.Range("E3").Value = "  -6.52%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "'321.12"
$ws.Range("E5").Value = "  +0.68%  "
$ws.Range("E6").Value = "  -9.48%  "
$ws.Range("D7").Value = "'0.578"
$ws.Range("E7").Value = "  -9.11%  "
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("D9").Value = "'0.562"
$ws.Range("E9").Value = "  -8.75%  "
$ws.Range("D10").Value = "'36.74"
$ws.Range("E10").Value = "  -10.53%  "
$ws.Range("D11").Value = "'54.01"
$ws.Range("E11").Value = "  -3.76%  "
$ws.Range("D12").Value = "'0.0826"
$ws.Range("E12").Value = "  -10.40%  "
$ws.Range("D13").Value = "'7.62"
$ws.Range("E13").Value = "  -10.74%  "
$ws.Range("E14").Value = "  -2.02%  "
$ws.Range("D15").Value = "'0.861"
$ws.Range("E15").Value = "  -12.51%  "
$ws.Range("D16").Value = "2.557.27"
$ws.Range("E16").Value = "  -6.37%  "
$ws.Range("D17").Value = "'14.34"
$ws.Range("E17").Value = "  -7.32%  "
$ws.Range("D18").Value = "2.218.19"
$ws.Range("E18").Value = "  -6.42%  "
$ws.Range("D19").Value = "42.892.73"
$ws.Range("E19").Value = "  -5.36%  "
$ws.Range("D20").Value = "'13.98"
$ws.Range("E20").Value = "  -9.06%  "
$ws.Range("E21").Value = "  -9.52%  "
$ws.Range("D22").Value = "'6.51"
$ws.Range("E22").Value = "  -11.20%  "
$ws.Range("B23").Value = "PancakeSwap"
$ws.Range("C23").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D23").Value = "'3.20"
$ws.Range("E23").Value = "  -11.73%  "
$ws.Range("B24").Value = "Litecoin"
$ws.Range("C24").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D24").Value = "'65.04"
$ws.Range("E24").Value = "  -11.33%  "
$ws.Range("D25").Value = "'236.27"
$ws.Range("E25").Value = "  -10.87%  "
$ws.Range("D26").Value = "'2.17"
$ws.Range("E26").Value = "  -7.20%  "
$ws.Range("E27").Value = "  -0.02%  "
$ws.Range("D28").Value = "'4.03"
$ws.Range("E28").Value = "  +1.19%  "
$ws.Range("D29").Value = "'9.96"
$ws.Range("E29").Value = "  -11.71%  "
$ws.Range("D30").Value = "'2.18"
$ws.Range("E30").Value = "  -5.58%  "
$ws.Range("E31").Value = "  -16.06%  "
$ws.Range("D32").Value = "'35.81"
$ws.Range("E32").Value = "  -4.07%  "
$ws.Range("D33").Value = "'20.27"
$ws.Range("E33").Value = "  -9.72%  "
$ws.Range("D34").Value = "'0.0862"
$ws.Range("E34").Value = "  -9.29%  "
$ws.Range("D35").Value = "'153.76"
$ws.Range("E35").Value = "  -8.96%  "
$ws.Range("E36").Value = "  -7.64%  "
$ws.Range("E37").Value = "  +2.90%  "
$ws.Range("E38").Value = "  -7.86%  "
$ws.Range("D39").Value = "'1.92"
$ws.Range("E39").Value = "  -0.70%  "
$ws.Range("E40").Value = "  -7.33%  "
$ws.Range("E41").Value = "  -11.66%  "
$ws.Range("D42").Value = "'3.66"
$ws.Range("E42").Value = "  -9.46%  "
$ws.Range("D43").Value = "'0.0321"
$ws.Range("E43").Value = "  -9.63%  "
$ws.Range("D44").Value = "'13.79"
$ws.Range("E44").Value = "  +5.99%  "
$ws.Range("E45").Value = "  -0.07%  "
$ws.Range("D46").Value = "1.722.75"
$ws.Range("E46").Value = "  -8.72%  "
$ws.Range("D47").Value = "'0.204"
$ws.Range("E47").Value = "  -10.87%  "
$ws.Range("D48").Value = "'84.59"
$ws.Range("E48").Value = "  -14.68%  "
$ws.Range("E49").Value = "  -12.78%  "
$ws.Range("D50").Value = "'8.79"
$ws.Range("E50").Value = "  -4.87%  "
$ws.Range("D51").Value = "'74.55"
$ws.Range("E51").Value = "  -11.98%  "
